# Fixed errors in excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpeedTest")

# Update measured/observed values in column O (the underlying raw data
# that drives the downstream formula columns T,U,V,W,X,Y via shared
# formulas). Excel will recalculate all dependents automatically.
$ws.Range("O7").Value = 1.7
$ws.Range("O8").Value = 1.7
$ws.Range("O9").Value = 1.8

$ws.Range("O12").Value = 3.4
$ws.Range("O13").Value = 3.4
$ws.Range("O14").Value = 3.5

$ws.Range("O21").Value = 3.2
$ws.Range("O22").Value = 3.2
$ws.Range("O23").Value = 3.2

$ws.Range("O34").Value = 2.8
$ws.Range("O35").Value = 2.8
$ws.Range("O36").Value = 3

$ws.Range("O40").Value = 4.8

# Update the frozen pane / selection state of the active sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("Y21").Select()

$excel.Calculate()
